# Auto-generated edit script: updates Hades_Profits workbook per scheduled runner diff.
# Applies numeric cell updates (and a few cell adds/removals) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value2 = 2333.3333
$ws.Range("J18").Value2 = 3000
$ws.Range("L18").Value2 = 3000
$ws.Range("N18").Value2 = -3568

$ws.Range("H40").Value2 = 1198.2858
$ws.Range("I40").Value2 = 1694
$ws.Range("J40").Value2 = 1000
$ws.Range("K40").Value2 = 1694
$ws.Range("L40").Value2 = 1000
$ws.Range("M40").Value2 = -1519
$ws.Range("N40").Value2 = -1350

$ws.Range("H125").Value2 = 805.4286
$ws.Range("I125").Value2 = 733.6667
$ws.Range("J125").Value2 = 1236
$ws.Range("K125").Value2 = 6603.0003
$ws.Range("L125").Value2 = 11124
$ws.Range("M125").Value2 = -4143.0003
$ws.Range("N125").Value2 = -16044

$ws.Range("H129").Value2 = 906.87933
$ws.Range("I129").Value2 = 700.17645
$ws.Range("J129").Value2 = 992.5854
$ws.Range("K129").Value2 = 2100.52935
$ws.Range("L129").Value2 = 2977.7562
$ws.Range("M129").Value2 = 2899.47065
$ws.Range("N129").Value2 = -12977.7562

$ws.Range("H137").Value2 = 2001457.8
$ws.Range("I137").Value2 = 7144164
$ws.Range("J137").Value2 = 1516.3055
$ws.Range("K137").Value2 = 21432492
$ws.Range("L137").Value2 = 4548.916499999999
$ws.Range("M137").Value2 = -21429942
$ws.Range("N137").Value2 = -9648.916499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value2 = 50000400
$ws.Range("I10").Value2 = 50000400
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 50000400
$ws.Range("L10").Value2 = 0
$ws.Range("M10").Value2 = -50000230
$ws.Range("N10").ClearContents()

$ws.Range("H32").Value2 = 19845.254
$ws.Range("I32").Value2 = 17346.865
$ws.Range("K32").Value2 = 17346.865
$ws.Range("M32").Value2 = -17059.865

$ws.Range("H61").Value2 = 95334500
$ws.Range("I61").Value2 = 66734460
$ws.Range("J61").Value2 = 166834600
$ws.Range("K61").Value2 = 66734460
$ws.Range("L61").Value2 = 166834600
$ws.Range("M61").Value2 = -66734248
$ws.Range("N61").Value2 = -166835024

$ws.Range("H74").Value2 = 18000774
$ws.Range("J74").Value2 = 500700
$ws.Range("L74").Value2 = 500700
$ws.Range("N74").Value2 = -502448

$ws.Range("H77").Value2 = 18000774
$ws.Range("J77").Value2 = 500700
$ws.Range("L77").Value2 = 2503500
$ws.Range("N77").Value2 = -2512236

$ws.Range("H97").Value2 = 3907906.2
$ws.Range("J97").Value2 = 905.5
$ws.Range("L97").Value2 = 905.5
$ws.Range("N97").Value2 = -1897.5

$ws.Range("H102").Value2 = 12988643
$ws.Range("I102").Value2 = 17858512
$ws.Range("K102").Value2 = 17858512
$ws.Range("M102").Value2 = -17856890

$ws.Range("H136").Value2 = 95334500
$ws.Range("I136").Value2 = 66734460
$ws.Range("J136").Value2 = 166834600
$ws.Range("K136").Value2 = 200203380
$ws.Range("L136").Value2 = 500503800
$ws.Range("M136").Value2 = -200200830
$ws.Range("N136").Value2 = -500508900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value2 = 1228.625
$ws.Range("I99").Value2 = 1191.6666
$ws.Range("K99").Value2 = 1191.6666
$ws.Range("M99").Value2 = 306.3334

$ws.Range("H134").Value2 = 6934
$ws.Range("J134").Value2 = 1800
$ws.Range("L134").Value2 = 5400
$ws.Range("N134").Value2 = -10470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H12").Value2 = 4242.5
$ws.Range("I12").Value2 = 985
$ws.Range("J12").Value2 = 7500
$ws.Range("K12").Value2 = 985
$ws.Range("L12").Value2 = 7500
$ws.Range("M12").Value2 = -815
$ws.Range("N12").Value2 = -7840

$ws.Range("H130").Value2 = 57282.105
$ws.Range("J130").Value2 = 57282.105
$ws.Range("L130").Value2 = 57282.105
$ws.Range("N130").Value2 = -67322.10500000001

$ws.Range("H135").Value2 = 48319.812
$ws.Range("J135").Value2 = 48319.812
$ws.Range("L135").Value2 = 48319.812
$ws.Range("N135").Value2 = -58459.812

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value2 = 0
$ws.Range("I20").Value2 = 0
$ws.Range("K20").Value2 = 0
$ws.Range("M20").ClearContents()

$ws.Range("H68").Value2 = 455734.38
$ws.Range("I68").Value2 = 914.6429000000001
$ws.Range("J68").Value2 = 1251669
$ws.Range("K68").Value2 = 2743.9287
$ws.Range("L68").Value2 = 3755007
$ws.Range("M68").Value2 = -1932.9287
$ws.Range("N68").Value2 = -3756629

$ws.Range("H71").Value2 = 455734.38
$ws.Range("I71").Value2 = 914.6429000000001
$ws.Range("J71").Value2 = 1251669
$ws.Range("K71").Value2 = 8231.786100000001
$ws.Range("L71").Value2 = 11265021
$ws.Range("M71").Value2 = -4175.786100000001
$ws.Range("N71").Value2 = -11273133

$ws.Range("H92").Value2 = 1115
$ws.Range("I92").Value2 = 830
$ws.Range("J92").Value2 = 1628
$ws.Range("K92").Value2 = 2490
$ws.Range("L92").Value2 = 4884
$ws.Range("M92").Value2 = -1242
$ws.Range("N92").Value2 = -7380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value2 = 2062.5264
$ws.Range("I102").Value2 = 1918.2667
$ws.Range("K102").Value2 = 1918.2667
$ws.Range("M102").Value2 = -296.2666999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 6850
$ws.Range("J40").Value2 = 6850
$ws.Range("L40").Value2 = 6850
$ws.Range("N40").Value2 = -7122

$ws.Range("H88").Value2 = 50555.555
$ws.Range("J88").Value2 = 50555.555
$ws.Range("L88").Value2 = 50555.555
$ws.Range("N88").Value2 = -51411.555

$ws.Range("H91").Value2 = 50555.555
$ws.Range("J91").Value2 = 50555.555
$ws.Range("L91").Value2 = 50555.555
$ws.Range("N91").Value2 = -53519.555

$ws.Range("H132").Value2 = 95407.63
$ws.Range("I132").Value2 = 1266.6666
$ws.Range("J132").Value2 = 130710.5
$ws.Range("K132").Value2 = 3799.9998
$ws.Range("L132").Value2 = 392131.5
$ws.Range("M132").Value2 = -1269.9998
$ws.Range("N132").Value2 = -397191.5

$ws.Range("H134").Value2 = 19493.77
$ws.Range("I134").Value2 = 8390
$ws.Range("J134").Value2 = 20419.084
$ws.Range("K134").Value2 = 8390
$ws.Range("L134").Value2 = 20419.084
$ws.Range("M134").Value2 = -3320
$ws.Range("N134").Value2 = -30559.084

$ws.Range("H135").Value2 = 53571.43
$ws.Range("J135").Value2 = 53571.43
$ws.Range("L135").Value2 = 53571.43
$ws.Range("N135").Value2 = -63711.43

$ws.Range("H136").Value2 = 183128.55
$ws.Range("I136").Value2 = 201122.2
$ws.Range("K136").Value2 = 603366.6000000001
$ws.Range("M136").Value2 = -600816.6000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value2 = 0
$ws.Range("J20").Value2 = 0
$ws.Range("L20").Value2 = 0
$ws.Range("N20").ClearContents()

$ws.Range("H46").Value2 = 50000
$ws.Range("J46").Value2 = 50000
$ws.Range("L46").Value2 = 50000
$ws.Range("N46").Value2 = -50462

$ws.Range("H100").Value2 = 100528.7
$ws.Range("I100").Value2 = 100760.4
$ws.Range("J100").Value2 = 100297
$ws.Range("K100").Value2 = 201520.8
$ws.Range("L100").Value2 = 200594
$ws.Range("M100").Value2 = -200979.8
$ws.Range("N100").Value2 = -201676

$ws.Range("H134").Value2 = 50000
$ws.Range("J134").Value2 = 50000
$ws.Range("L134").Value2 = 150000
$ws.Range("N134").Value2 = -155070

$ws.Range("H135").Value2 = 44333
$ws.Range("J135").Value2 = 44333
$ws.Range("L135").Value2 = 44333
$ws.Range("N135").Value2 = -54473

$ws.Range("H136").Value2 = 49461.977
$ws.Range("I136").Value2 = 32429.469
$ws.Range("K136").Value2 = 97288.40700000001
$ws.Range("M136").Value2 = -94738.40700000001

